$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Right" count for Marking row (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update "Right" count for Total row (B12): 51 -> 85
$ws.Range("B12").Value = 85

# Update Correct/Total marks string (E12): "49/84" -> "85/140"
$ws.Range("E12").Value = "85/140"
